$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Forces a value to be stored as text even if it looks like a number
    # or a date, without leaving a visible style/format on the cell.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# --- Rows 22 and 23 swap their species-observation data -------------------
# (keep P/S/T/U/V/W/Y/AA/AD/AE/AG/AT/AW/AX/AY as-is; they are identical
#  between the two rows already)

$ws.Range("A22").Value = 130834063
$ws.Range("B22").Value = 91771
Set-TextValue $ws.Range("D22") "LC"
$ws.Range("E22").Value = 5447
Set-TextValue $ws.Range("F22") "Vedticka"
Set-TextValue $ws.Range("G22") "Fuscoporia viticola"
Set-TextValue $ws.Range("H22") "(Schwein.) Murrill"
$ws.Range("I22").ClearContents()
$ws.Range("Q22").Value = 516858
$ws.Range("R22").Value = 6677494

$ws.Range("A23").Value = 130834071
$ws.Range("B23").Value = 57897
Set-TextValue $ws.Range("D23") "NT"
$ws.Range("E23").Value = 100048
Set-TextValue $ws.Range("F23") "Mindre hackspett"
Set-TextValue $ws.Range("G23") "Dryobates minor"
Set-TextValue $ws.Range("H23") "(Linnaeus, 1758)"
Set-TextValue $ws.Range("I23") "1"
$ws.Range("Q23").Value = 516848
$ws.Range("R23").Value = 6677496

# --- Rows 26 and 27 swap their species-observation data -------------------
# (D/H/I/P/S/T/U/V/W/AD/AE/AG/AT/AW/AX/AY are identical between the two
#  rows already, so they are left untouched)

$ws.Range("A26").Value = 130834054
$ws.Range("B26").Value = 97878
$ws.Range("E26").Value = 221945
Set-TextValue $ws.Range("F26") "Revlummer"
Set-TextValue $ws.Range("G26") "Lycopodium annotinum"
$ws.Range("Q26").Value = 516969
$ws.Range("R26").Value = 6677335
Set-TextValue $ws.Range("Y26") "2025-09-25"
Set-TextValue $ws.Range("AA26") "2025-09-25"

$ws.Range("A27").Value = 130834060
$ws.Range("B27").Value = 97881
$ws.Range("E27").Value = 221946
Set-TextValue $ws.Range("F27") "Mattlummer"
Set-TextValue $ws.Range("G27") "Lycopodium clavatum"
$ws.Range("Q27").Value = 516909
$ws.Range("R27").Value = 6677366
Set-TextValue $ws.Range("Y27") "2025-09-24"
Set-TextValue $ws.Range("AA27") "2025-09-24"
